$d = $word.ActiveDocument

# --- Table of contents: "Page des X" -> "Onglet des X" -------------------
$d.Content.Find.Execute("Page des pays", $true, $false, $false, $false, $false, $true, 1, $false, "Onglet des pays", 2)
$d.Content.Find.Execute("Page des départements", $true, $false, $false, $false, $false, $true, 1, $false, "Onglet des départements", 2)
$d.Content.Find.Execute("Page des médecins", $true, $false, $false, $false, $false, $true, 1, $false, "Onglet des médecins", 2)

# --- Table of contents: "Modifier un X" -> "Modifier/Supprimer un X" -----
$d.Content.Find.Execute(".   Modifier un pays", $true, $false, $false, $false, $false, $true, 1, $false, ".   Modifier/Supprimer un pays", 2)
$d.Content.Find.Execute(".   Modifier un département", $true, $false, $false, $false, $false, $true, 1, $false, ".   Modifier/Supprimer un département", 2)
$d.Content.Find.Execute(".   Modifier un médecin", $true, $false, $false, $false, $false, $true, 1, $false, ".   Modifier/Supprimer un médecin", 2)

# --- Table of contents: merge split runs for "Ajouter un département" ----
$d.Content.Find.Execute(".   Ajouter un département", $true, $false, $false, $false, $false, $true, 1, $false, ".   Ajouter un département", 2)

# --- Footer 2: merge "GSB " + "Médecins...Page" runs ----------------------
$d.Content.Find.Execute("GSB Médecins                                                                                                                      Page", $true, $false, $false, $false, $false, $true, 1, $false, "GSB Médecins                                                                                                                      Page", 2)

# --- Header 1: merge "Documentation utilisateur réalisée le " + "20 mars " + "2021" ----
$d.Content.Find.Execute("Documentation utilisateur réalisée le 20 mars 2021", $true, $false, $false, $false, $false, $true, 1, $false, "Documentation utilisateur réalisée le 20 mars 2021", 2)
